# Applies numeric value updates to the Leve profit-tracking tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# as produced by the scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3350.2454
$ws.Range("I76").Value = 2947.5642
$ws.Range("J76").Value = 4472
$ws.Range("K76").Value = 2947.5642
$ws.Range("L76").Value = 4472
$ws.Range("M76").Value = -2632.5642
$ws.Range("N76").Value = -5102
$ws.Range("H79").Value = 3350.2454
$ws.Range("I79").Value = 2947.5642
$ws.Range("J79").Value = 4472
$ws.Range("K79").Value = 2947.5642
$ws.Range("L79").Value = 4472
$ws.Range("M79").Value = -1855.5642
$ws.Range("N79").Value = -6656
$ws.Range("H86").Value = 43785.645
$ws.Range("I86").Value = 23077.666
$ws.Range("J86").Value = 81060
$ws.Range("K86").Value = 23077.666
$ws.Range("L86").Value = 81060
$ws.Range("M86").Value = -21954.666
$ws.Range("N86").Value = -83306
$ws.Range("H89").Value = 43785.645
$ws.Range("I89").Value = 23077.666
$ws.Range("J89").Value = 81060
$ws.Range("K89").Value = 115388.33
$ws.Range("L89").Value = 405300
$ws.Range("M89").Value = -109772.33
$ws.Range("N89").Value = -416532
$ws.Range("H125").Value = 1977995.6
$ws.Range("I125").Value = 33866.668
$ws.Range("J125").Value = 3922124.8
$ws.Range("K125").Value = 304800.012
$ws.Range("L125").Value = 35299123.2
$ws.Range("M125").Value = -302340.012
$ws.Range("N125").Value = -35304043.2
$ws.Range("H129").Value = 963.5862
$ws.Range("I129").Value = 448.6
$ws.Range("J129").Value = 1070.875
$ws.Range("K129").Value = 1345.8
$ws.Range("L129").Value = 3212.625
$ws.Range("M129").Value = 3654.2
$ws.Range("N129").Value = -13212.625
$ws.Range("H132").Value = 9702.387000000001
$ws.Range("I132").Value = 14307.728
$ws.Range("J132").Value = 5097.0454
$ws.Range("K132").Value = 42923.18399999999
$ws.Range("L132").Value = 15291.1362
$ws.Range("M132").Value = -40393.18399999999
$ws.Range("N132").Value = -20351.1362
$ws.Range("H138").Value = 3354.73
$ws.Range("I138").Value = 1837.3422
$ws.Range("J138").Value = 4284.7417
$ws.Range("K138").Value = 5512.0266
$ws.Range("L138").Value = 12854.2251
$ws.Range("M138").Value = -372.0266000000001
$ws.Range("N138").Value = -23134.2251
$ws.Range("H141").Value = 1013.03705
$ws.Range("I141").Value = 1013.03705
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3039.11115
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2140.88885
$ws.Range("N141").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2261.5
$ws.Range("I2").Value = 2261.5
$ws.Range("K2").Value = 2261.5
$ws.Range("M2").Value = -2148.5
$ws.Range("H32").Value = 2353.77
$ws.Range("I32").Value = 2332.6082
$ws.Range("J32").Value = 3038
$ws.Range("K32").Value = 2332.6082
$ws.Range("L32").Value = 3038
$ws.Range("M32").Value = -2045.6082
$ws.Range("N32").Value = -3612
$ws.Range("H116").Value = 2261.5
$ws.Range("I116").Value = 2261.5
$ws.Range("K116").Value = 2261.5
$ws.Range("M116").Value = 32.5
$ws.Range("H132").Value = 1771.9166
$ws.Range("I132").Value = 1468.9143
$ws.Range("J132").Value = 2587.6924
$ws.Range("K132").Value = 4406.742899999999
$ws.Range("L132").Value = 7763.0772
$ws.Range("M132").Value = -1876.742899999999
$ws.Range("N132").Value = -12823.0772

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2261.5
$ws.Range("I3").Value = 2261.5
$ws.Range("K3").Value = 2261.5
$ws.Range("M3").Value = -2147.5
$ws.Range("H51").Value = 36590
$ws.Range("J51").Value = 36590
$ws.Range("L51").Value = 36590
$ws.Range("N51").Value = -37572
$ws.Range("H134").Value = 1250.963
$ws.Range("I134").Value = 905.8421
$ws.Range("J134").Value = 2070.625
$ws.Range("K134").Value = 2717.5263
$ws.Range("L134").Value = 6211.875
$ws.Range("M134").Value = -182.5263
$ws.Range("N134").Value = -11281.875

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 40035.5
$ws.Range("J47").Value = 40035.5
$ws.Range("L47").Value = 40035.5
$ws.Range("N47").Value = -41167.5
$ws.Range("H132").Value = 963705.25
$ws.Range("I132").Value = 1924601.8
$ws.Range("J132").Value = 2808.7693
$ws.Range("K132").Value = 5773805.4
$ws.Range("L132").Value = 8426.3079
$ws.Range("M132").Value = -5771275.4
$ws.Range("N132").Value = -13486.3079
$ws.Range("H134").Value = 2245.9375
$ws.Range("I134").Value = 2446.7827
$ws.Range("J134").Value = 1732.6666
$ws.Range("K134").Value = 7340.348100000001
$ws.Range("L134").Value = 5197.9998
$ws.Range("M134").Value = -4805.348100000001
$ws.Range("N134").Value = -10267.9998
$ws.Range("H141").Value = 40354.688
$ws.Range("I141").Value = 28000
$ws.Range("J141").Value = 42119.645
$ws.Range("K141").Value = 28000
$ws.Range("L141").Value = 42119.645
$ws.Range("M141").Value = -22820
$ws.Range("N141").Value = -52479.645

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 666.6667
$ws.Range("I60").Value = 600
$ws.Range("J60").Value = 800
$ws.Range("K60").Value = 1800
$ws.Range("L60").Value = 2400
$ws.Range("M60").Value = -1549
$ws.Range("N60").Value = -2902
$ws.Range("H129").Value = 1471.4872
$ws.Range("I129").Value = 945.55554
$ws.Range("J129").Value = 1629.2667
$ws.Range("K129").Value = 2836.66662
$ws.Range("L129").Value = 4887.800099999999
$ws.Range("M129").Value = 2163.33338
$ws.Range("N129").Value = -14887.8001
$ws.Range("H131").Value = 841.96
$ws.Range("J131").Value = 854.1429000000001
$ws.Range("L131").Value = 2562.4287
$ws.Range("N131").Value = -12642.4287

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2960
$ws.Range("I80").Value = 2940
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2940
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1942
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2960
$ws.Range("I83").Value = 2940
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 14700
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -9708
$ws.Range("N83").Value = -24984
$ws.Range("H132").Value = 2847.2593
$ws.Range("I132").Value = 2440.2856
$ws.Range("J132").Value = 3285.5386
$ws.Range("K132").Value = 7320.8568
$ws.Range("L132").Value = 9856.6158
$ws.Range("M132").Value = -4790.8568
$ws.Range("N132").Value = -14916.6158
$ws.Range("H138").Value = 43388.645
$ws.Range("J138").Value = 43388.645
$ws.Range("L138").Value = 43388.645
$ws.Range("N138").Value = -53668.645

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 70564.87
$ws.Range("I7").Value = 94497.63
$ws.Range("J7").Value = 4749.75
$ws.Range("K7").Value = 94497.63
$ws.Range("L7").Value = 4749.75
$ws.Range("M7").Value = -94385.63
$ws.Range("N7").Value = -4973.75
$ws.Range("H40").Value = 27604.5
$ws.Range("I40").Value = 37663.57
$ws.Range("J40").Value = 4133.3335
$ws.Range("K40").Value = 37663.57
$ws.Range("L40").Value = 4133.3335
$ws.Range("M40").Value = -37527.57
$ws.Range("N40").Value = -4405.3335
$ws.Range("H45").Value = 33000
$ws.Range("J45").Value = 33000
$ws.Range("L45").Value = 33000
$ws.Range("N45").Value = -33814
$ws.Range("H93").Value = 929
$ws.Range("I93").Value = 923.75
$ws.Range("K93").Value = 923.75
$ws.Range("M93").Value = 324.25
$ws.Range("H126").Value = 70564.87
$ws.Range("I126").Value = 94497.63
$ws.Range("J126").Value = 4749.75
$ws.Range("K126").Value = 283492.89
$ws.Range("L126").Value = 14249.25
$ws.Range("M126").Value = -281022.89
$ws.Range("N126").Value = -19189.25
$ws.Range("H127").Value = 53235
$ws.Range("J127").Value = 53235
$ws.Range("L127").Value = 53235
$ws.Range("N127").Value = -63155
$ws.Range("H132").Value = 6076.194
$ws.Range("I132").Value = 6276.638
$ws.Range("J132").Value = 5605.15
$ws.Range("K132").Value = 18829.914
$ws.Range("L132").Value = 16815.45
$ws.Range("M132").Value = -16299.914
$ws.Range("N132").Value = -21875.45

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 42485.734
$ws.Range("J46").Value = 42485.734
$ws.Range("L46").Value = 42485.734
$ws.Range("N46").Value = -42947.734
$ws.Range("H96").Value = 1684
$ws.Range("J96").Value = 1684
$ws.Range("L96").Value = 1684
$ws.Range("N96").Value = -4430
$ws.Range("H113").Value = 1037.25
$ws.Range("I113").Value = 1278.1666
$ws.Range("J113").Value = 314.5
$ws.Range("K113").Value = 3834.4998
$ws.Range("L113").Value = 943.5
$ws.Range("M113").Value = -1664.4998
$ws.Range("N113").Value = -5283.5
$ws.Range("H134").Value = 42485.734
$ws.Range("J134").Value = 42485.734
$ws.Range("L134").Value = 127457.202
$ws.Range("N134").Value = -132527.202
$ws.Range("H136").Value = 2647.3877
$ws.Range("I136").Value = 574.4091
$ws.Range("J136").Value = 4336.4814
$ws.Range("K136").Value = 1723.2273
$ws.Range("L136").Value = 13009.4442
$ws.Range("M136").Value = 826.7727
$ws.Range("N136").Value = -18109.4442

Write-Host "Updated 234 cells across 8 sheets."